$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Rename the "View" field header to "Cache"
$ws.Range("F1").Value = "Cache"

# Set default value of the field to FALSE for all data rows
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
}
